# Daily attendance processing - 2025-10-09 14:20:28
#
# Column G on the "Session Analysis Results" sheet holds a comma-separated
# list of the users/processes that touched a given attendance session
# (e.g. "dnasr281@gmail.com, System"). The recorder that appends "System"
# to the end of that list should instead be putting it first. This pass
# walks every row of the sheet and, wherever the last entry in the list is
# "System"/"system" and it isn't already the first entry, swaps the first
# and last entries so "System" leads the list (any entries in between are
# left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ", ") { continue }

    $tokens = $val -split ", "
    $count = $tokens.Length
    if ($count -lt 2) { continue }

    $first = $tokens[0]
    $last = $tokens[$count - 1]

    # Case-sensitive compare (the -ceq/-cne operators in this host behave
    # case-insensitively, so fall back to the .NET String.Equals overload,
    # which is ordinal): "System" already in front should be left alone,
    # but a lowercase "system" in front still needs to be swapped to the end
    # (the pair relabels which occurrence is capitalised).
    if (($last.ToLower() -eq "system") -and (-not $first.Equals("System"))) {
        if ($count -gt 2) {
            $middle = $tokens[1..($count - 2)]
        } else {
            $middle = @()
        }
        $newTokens = @($last) + @($middle) + @($first)
        $newVal = $newTokens -join ", "
        $cell.Value = $newVal
    }
}
